$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.142.18'
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").Value = '2.371.91'
$ws.Range("E3").Value = '  +1.53%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''302.79'
$ws.Range("E5").Value = '  -0.23%  '

$ws.Range("D6").Value = '''95.48'
$ws.Range("E6").Value = '  +1.70%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '''0.503'
$ws.Range("E7").Value = '  +0.46%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("D9").Value = '''0.482'
$ws.Range("E9").Value = '  -2.50%  '

$ws.Range("D10").Value = '''34.39'
$ws.Range("E10").Value = '  +1.14%  '

$ws.Range("E11").Value = '  +3.72%  '

$ws.Range("D12").Value = '''0.0787'
$ws.Range("E12").Value = '  +0.65%  '

$ws.Range("D13").Value = '''18.29'
$ws.Range("E13").Value = '  -2.09%  '

$ws.Range("D14").Value = '''6.76'
$ws.Range("E14").Value = '  +0.80%  '

$ws.Range("D15").Value = '2.736.33'
$ws.Range("E15").Value = '  +1.28%  '

$ws.Range("D16").Value = '2.396.18'
$ws.Range("E16").Value = '  +0.98%  '

$ws.Range("D17").Value = '''0.799'
$ws.Range("E17").Value = '  +1.00%  '

$ws.Range("D18").Value = '43.122.77'
$ws.Range("E18").Value = '  +0.67%  '

$ws.Range("D19").Value = '''11.95'
$ws.Range("E19").Value = '  -0.82%  '

$ws.Range("D20").Value = '''6.27'
$ws.Range("E20").Value = '  +1.09%  '

$ws.Range("D21").Value = '0.0₃0888'
$ws.Range("E21").Value = '  +0.11%  '

$ws.Range("D22").Value = '''67.94'
$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("D23").Value = '''235.88'
$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("D24").Value = '''2.20'
$ws.Range("E24").Value = '  -0.32%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = '''2.42'
$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("D27").Value = '''24.50'
$ws.Range("E27").Value = '  -0.16%  '

$ws.Range("E28").Value = '  +15.00%  '

$ws.Range("D29").Value = '''9.36'
$ws.Range("E29").Value = '  +2.61%  '

$ws.Range("D30").Value = '''32.09'
$ws.Range("E30").Value = '  +2.58%  '

$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("D32").Value = '''5.02'
$ws.Range("E32").Value = '  +0.97%  '

$ws.Range("D33").Value = '''17.59'
$ws.Range("E33").Value = '  +2.25%  '

$ws.Range("E34").Value = '  +8.34%  '

$ws.Range("D35").Value = '''0.0730'
$ws.Range("E35").Value = '  -1.90%  '

$ws.Range("D36").Value = '''127.37'
$ws.Range("E36").Value = '  +1.97%  '

$ws.Range("E37").Value = '  +1.30%  '

$ws.Range("D38").Value = '''2.86'
$ws.Range("E38").Value = '  +3.67%  '

$ws.Range("D39").Value = '''4.31'
$ws.Range("E39").Value = '  -1.24%  '

$ws.Range("D40").Value = '''2.26'
$ws.Range("E40").Value = '  -2.41%  '

$ws.Range("E41").Value = '  -0.28%  '

$ws.Range("D42").Value = '''20.80'
$ws.Range("E42").Value = '  -6.95%  '

$ws.Range("D43").Value = '1.930.54'
$ws.Range("E43").Value = '  -0.24%  '

$ws.Range("E44").Value = '  -1.01%  '

$ws.Range("E45").Value = '  +2.50%  '

$ws.Range("D46").Value = '''9.22'
$ws.Range("E46").Value = '  -9.27%  '

$ws.Range("D47").Value = '''2.73'
$ws.Range("E47").Value = '  +0.99%  '

$ws.Range("D48").Value = '2.597.13'
$ws.Range("E48").Value = '  +1.15%  '

$ws.Range("D49").Value = '''1.51'
$ws.Range("E49").Value = '  +2.85%  '

$ws.Range("D50").Value = '''71.42'
$ws.Range("E50").Value = '  -0.08%  '

$ws.Range("D51").Value = '''51.53'
$ws.Range("E51").Value = '  -2.33%  '
